$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (percentages) and two-dot "thousands" style column D values ---
# These new strings are not auto-recognized by Excel as numbers, so a plain
# Value assignment keeps them as text (inline/shared string).
$ws.Range('D2').Value = '26.193.96'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '1.676.72'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('E4').Value = '  -0.71%  '
$ws.Range('E5').Value = '  -3.53%  '
$ws.Range('E6').Value = '  -4.47%  '
$ws.Range('E7').Value = '  -0.70%  '
$ws.Range('E8').Value = '  -3.18%  '
$ws.Range('E9').Value = '  -2.86%  '
$ws.Range('E10').Value = '  -3.22%  '
$ws.Range('E11').Value = '  -1.90%  '
$ws.Range('D12').Value = '1.674.71'
$ws.Range('E12').Value = '  -1.64%  '
$ws.Range('E13').Value = '  -2.32%  '
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('E16').Value = '  -4.70%  '
$ws.Range('D17').Value = '26.259.78'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('E19').Value = '  -2.91%  '
$ws.Range('E20').Value = '  -2.33%  '
$ws.Range('E21').Value = '  -5.08%  '
$ws.Range('E22').Value = '  -1.47%  '
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('E25').Value = '  -5.98%  '
$ws.Range('E26').Value = '  -4.56%  '
$ws.Range('E27').Value = '  +1.22%  '
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('E30').Value = '  -3.82%  '
$ws.Range('E31').Value = '  -3.62%  '
$ws.Range('E32').Value = '  -4.73%  '
$ws.Range('E33').Value = '  -3.57%  '
$ws.Range('E34').Value = '  -4.26%  '
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('E37').Value = '  -0.85%  '
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('E39').Value = '  -2.47%  '
$ws.Range('D40').Value = '1.085.24'
$ws.Range('E40').Value = '  -3.38%  '
$ws.Range('E41').Value = '  -1.78%  '
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('E43').Value = '  -1.49%  '
$ws.Range('D44').Value = '1.828.52'
$ws.Range('E44').Value = '  -1.45%  '
$ws.Range('E45').Value = '  +2.61%  '
$ws.Range('E46').Value = '  -2.86%  '
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('E48').Value = '  -2.98%  '
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('E51').Value = '  -2.72%  '

# --- Column D price values that look like plain numbers ---
# Force a text number format first so Excel does not coerce these into
# numeric values, then restore the default "Normal" style so no extra
# explicit cell style is left behind.
$numericPriceCells = @('D5', 'D6', 'D9', 'D11', 'D13', 'D14', 'D15', 'D16', 'D20', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D29', 'D31', 'D32', 'D33', 'D34', 'D35', 'D37', 'D38', 'D41', 'D43', 'D47', 'D48', 'D49', 'D51')
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D5').Value = '211.79'
$ws.Range('D6').Value = '0.5253'
$ws.Range('D9').Value = '0.06295'
$ws.Range('D11').Value = '0.07558'
$ws.Range('D13').Value = '4.455'
$ws.Range('D14').Value = '0.5635'
$ws.Range('D15').Value = '67.02'
$ws.Range('D16').Value = '0.000008024'
$ws.Range('D20').Value = '188.07'
$ws.Range('D21').Value = '10.44'
$ws.Range('D22').Value = '6.191'
$ws.Range('D24').Value = '149.50'
$ws.Range('D25').Value = '0.1250'
$ws.Range('D26').Value = '7.577'
$ws.Range('D27').Value = '16.06'
$ws.Range('D29').Value = '1.358'
$ws.Range('D31').Value = '3.499'
$ws.Range('D32').Value = '3.441'
$ws.Range('D33').Value = '1.635'
$ws.Range('D34').Value = '1.003'
$ws.Range('D35').Value = '0.6071'
$ws.Range('D37').Value = '2.743'
$ws.Range('D38').Value = '6.100'
$ws.Range('D41').Value = '0.8693'
$ws.Range('D43').Value = '100.08'
$ws.Range('D47').Value = '0.9979'
$ws.Range('D48').Value = '8.003'
$ws.Range('D49').Value = '0.05234'
$ws.Range('D51').Value = '5.980'

foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).Style = "Normal"
}
